# Weekly update: insert two new price observations at the top of the
# "Ramas de apio" (Vega Modelo de Temuco) data range, pushing older rows
# down (mirrors how this logica_diaria workbook accretes one row per
# market visit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow($rowNum, $fecha, $volumen, $precioMin, $precioMax, $precioProm, $precioKg) {
    $ws.Cells.Item($rowNum, 1).Value = 10
    $ws.Cells.Item($rowNum, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($rowNum, 3).Value = "La Araucanía"
    $ws.Cells.Item($rowNum, 4).Value = $fecha
    $ws.Cells.Item($rowNum, 5).Value = 9
    $ws.Cells.Item($rowNum, 6).Value = 100112017
    $ws.Cells.Item($rowNum, 7).Value = "Ramas de apio"
    $ws.Cells.Item($rowNum, 8).Value = "Sin especificar"
    $ws.Cells.Item($rowNum, 9).Value = "Primera"
    $ws.Cells.Item($rowNum, 10).Value = $volumen
    $ws.Cells.Item($rowNum, 11).Value = $precioMin
    $ws.Cells.Item($rowNum, 12).Value = $precioMax
    $ws.Cells.Item($rowNum, 13).Value = $precioProm
    $ws.Cells.Item($rowNum, 14).Value = "`$/paquete"
    $ws.Cells.Item($rowNum, 15).Value = "Región de La Araucanía"
    $ws.Cells.Item($rowNum, 16).Value = $precioKg
    $ws.Cells.Item($rowNum, 17).Value = 1
    $ws.Cells.Item($rowNum, 18).Value = "Hortaliza"
}

# Insert a brand new row at row 10, shifting the existing row 10 (and
# everything below) down by one.
$ws.Rows.Item(10).Insert()
Set-DataRow 10 44679 50 5000 5000 5000 5000

# Insert a second new row after the row that now holds the former row 10
# (i.e. at row 12), shifting everything from the former row 11 onward
# down by one more.
$ws.Rows.Item(12).Insert()
Set-DataRow 12 44680 20 5000 5000 5000 5000
